$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1954397394136808
$ws.Range("C2").Value = 0.5504885993485342
$ws.Range("J2").Value = 0.01628664495114007
$ws.Range("P2").Value = 0.1465798045602606
$ws.Range("S2").Value = 0.09120521172638436
$ws.Range("B3").Value = 0.01694915254237288
$ws.Range("C3").Value = 0.03389830508474576
$ws.Range("J3").Value = 0.03389830508474576
$ws.Range("P3").Value = 0.7570621468926554
$ws.Range("S3").Value = 0.1581920903954802
$ws.Range("J4").Value = 0.04545454545454546
$ws.Range("P4").Value = 0.6136363636363636
$ws.Range("S4").Value = 0.3409090909090909
$ws.Range("P5").Value = 0.75
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.05741626794258373
$ws.Range("F6").Value = 0.04784688995215311
$ws.Range("J6").Value = 0.2440191387559809
$ws.Range("O6").Value = 0.009569377990430622
$ws.Range("Q6").Value = 0.1961722488038277
$ws.Range("R6").Value = 0.06698564593301436
$ws.Range("S6").Value = 0.3779904306220095
$ws.Range("B7").Value = 0.1267605633802817
$ws.Range("D7").Value = 0.02347417840375587
$ws.Range("F7").Value = 0.04694835680751173
$ws.Range("J7").Value = 0.1032863849765258
$ws.Range("O7").Value = 0.04225352112676056
$ws.Range("Q7").Value = 0.1971830985915493
$ws.Range("R7").Value = 0.04225352112676056
$ws.Range("S7").Value = 0.4178403755868544
$ws.Range("B8").Value = 0.09380863039399624
$ws.Range("D8").Value = 0.01688555347091933
$ws.Range("E8").Value = 0.00375234521575985
$ws.Range("F8").Value = 0.05628517823639775
$ws.Range("J8").Value = 0.1200750469043152
$ws.Range("O8").Value = 0.01688555347091933
$ws.Range("Q8").Value = 0.1538461538461539
$ws.Range("R8").Value = 0.06566604127579738
$ws.Range("S8").Value = 0.4727954971857411
$ws.Range("B9").Value = 0.09049773755656108
$ws.Range("D9").Value = 0.01809954751131222
$ws.Range("F9").Value = 0.05882352941176471
$ws.Range("J9").Value = 0.09954751131221719
$ws.Range("O9").Value = 0.03167420814479638
$ws.Range("Q9").Value = 0.1583710407239819
$ws.Range("R9").Value = 0.05882352941176471
$ws.Range("S9").Value = 0.4841628959276018
$ws.Range("B10").Value = 0.1090604026845638
$ws.Range("D10").Value = 0.02348993288590604
$ws.Range("E10").Value = 0.001677852348993289
$ws.Range("F10").Value = 0.06124161073825504
$ws.Range("J10").Value = 0.1107382550335571
$ws.Range("O10").Value = 0.01677852348993289
$ws.Range("Q10").Value = 0.1971476510067114
$ws.Range("R10").Value = 0.06963087248322147
$ws.Range("S10").Value = 0.410234899328859
$ws.Range("G11").Value = 0.1861861861861862
$ws.Range("J11").Value = 0.06606606606606606
$ws.Range("K11").Value = 0.2432432432432433
$ws.Range("L11").Value = 0.4984984984984985
$ws.Range("S11").Value = 0.006006006006006006
$ws.Range("G12").Value = 0.7485714285714286
$ws.Range("J12").Value = 0.2057142857142857
$ws.Range("K12").Value = 0.005714285714285714
$ws.Range("L12").Value = 0.03428571428571429
$ws.Range("S12").Value = 0.005714285714285714
$ws.Range("G13").Value = 0.5869565217391305
$ws.Range("J13").Value = 0.3260869565217391
$ws.Range("S13").Value = 0.08695652173913043
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.1666666666666667
$ws.Range("S14").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.02105263157894737
$ws.Range("H15").Value = 0.1157894736842105
$ws.Range("I15").Value = 0.07368421052631578
$ws.Range("J15").Value = 0.2842105263157895
$ws.Range("K15").Value = 0.07894736842105263
$ws.Range("M15").Value = 0.01052631578947368
$ws.Range("N15").Value = 0.005263157894736842
$ws.Range("O15").Value = 0.06315789473684211
$ws.Range("S15").Value = 0.3473684210526316
$ws.Range("H16").Value = 0.2572815533980582
$ws.Range("I16").Value = 0.08737864077669903
$ws.Range("J16").Value = 0.3349514563106796
$ws.Range("K16").Value = 0.0970873786407767
$ws.Range("M16").Value = 0.01456310679611651
$ws.Range("N16").Value = 0.004854368932038835
$ws.Range("O16").Value = 0.04854368932038835
$ws.Range("S16").Value = 0.1553398058252427
$ws.Range("F17").Value = 0.02073732718894009
$ws.Range("H17").Value = 0.1866359447004608
$ws.Range("I17").Value = 0.1013824884792627
$ws.Range("J17").Value = 0.4032258064516129
$ws.Range("K17").Value = 0.1152073732718894
$ws.Range("M17").Value = 0.01612903225806452
$ws.Range("N17").Value = 0.002304147465437788
$ws.Range("O17").Value = 0.04377880184331797
$ws.Range("S17").Value = 0.1105990783410138
$ws.Range("F18").Value = 0.04545454545454546
$ws.Range("H18").Value = 0.1688311688311688
$ws.Range("I18").Value = 0.1233766233766234
$ws.Range("J18").Value = 0.3896103896103896
$ws.Range("K18").Value = 0.09740259740259741
$ws.Range("O18").Value = 0.03246753246753246
$ws.Range("S18").Value = 0.1428571428571428
$ws.Range("F19").Value = 0.0173973556019485
$ws.Range("H19").Value = 0.2456506610995129
$ws.Range("I19").Value = 0.08768267223382047
$ws.Range("J19").Value = 0.3298538622129436
$ws.Range("K19").Value = 0.1064718162839248
$ws.Range("M19").Value = 0.02644398051496172
$ws.Range("N19").Value = 0.002783576896311761
$ws.Range("O19").Value = 0.05149617258176757
$ws.Range("S19").Value = 0.1322199025748086
